# Update "想去人数" (want-to-go count) values in column F for both the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets, per gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - first worksheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1269
$ws1.Range("F13").Value = 70
$ws1.Range("F14").Value = 44
$ws1.Range("F15").Value = 57
$ws1.Range("F16").Value = 604
$ws1.Range("F18").Value = 751
$ws1.Range("F22").Value = 60
$ws1.Range("F23").Value = 66
$ws1.Range("F24").Value = 2668
$ws1.Range("F25").Value = 5174
$ws1.Range("F29").Value = 3075
$ws1.Range("F31").Value = 2255
$ws1.Range("F35").Value = 123
$ws1.Range("F38").Value = 25
$ws1.Range("F39").Value = 462
$ws1.Range("F41").Value = 30
$ws1.Range("F44").Value = 39

# Sheet "全部类型" - fourth worksheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1269
$ws4.Range("F13").Value = 70
$ws4.Range("F15").Value = 44
$ws4.Range("F16").Value = 57
$ws4.Range("F17").Value = 604
$ws4.Range("F19").Value = 751
$ws4.Range("F23").Value = 60
$ws4.Range("F24").Value = 66
$ws4.Range("F25").Value = 2668
$ws4.Range("F26").Value = 5174
$ws4.Range("F30").Value = 3075
$ws4.Range("F32").Value = 2255
$ws4.Range("F36").Value = 123
$ws4.Range("F39").Value = 25
$ws4.Range("F40").Value = 462
$ws4.Range("F42").Value = 30
$ws4.Range("F45").Value = 39
